# Adds 2 new studies (4 rows of data) to the LDL Study Summary sheet,
# mirroring the commit "Added 2 studies to summary excel sheet".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oupLink  = "https://academic.oup.com/jcem/article/88/4/1617/2845298?login=false"
$nejmLink = "https://www.nejm.org/doi/10.1056/NEJMoa022207?url_ver=Z39.88-2003&rfr_id=ori:rid:crossref.org&rfr_dat=cr_pub%20%200www.ncbi.nlm.nih.gov"

# ---- Row 31: DAlessio2003, 3m ----
$ws.Range("A31").Value = "DAlessio2003"
$ws.Range("B31").Value = 91.2
$ws.Range("C31").Value = -7.6
$ws.Range("D31").Value = 33.17
$ws.Range("F31").Value = 124.86
$ws.Range("G31").Value = -11.86
$ws.Range("H31").Value = 5.39
$ws.Range("I31").Value = 5.34
$ws.Range("L31").Value = 22
$ws.Range("M31").Value = "3m"
$ws.Range("S31").Value = $oupLink
$ws.Range("P31").Value = "CHO too high"

# ---- Row 32: DAlessio2003, 6m (with real hyperlink on S32) ----
$ws.Range("A32").Value = "DAlessio2003"
$ws.Range("B32").Value = 91.2
$ws.Range("C32").Value = -8.5
$ws.Range("D32").Value = 33.17
$ws.Range("F32").Value = 124.86
$ws.Range("G32").Value = -0.86
$ws.Range("H32").Value = 5.39
$ws.Range("I32").Value = 5.81
$ws.Range("L32").Value = 22
$ws.Range("M32").Value = "6m"
$ws.Range("S32").Value = $oupLink
$null = $ws.Hyperlinks.Add($ws.Range("S32"), $oupLink)
$ws.Range("S32").Style = "Hyperlink"
$ws.Range("P32").Value = "CHO too high"

# ---- Row 33: Klein2003, 3m ----
$ws.Range("A33").Value = "Klein2003"

# ---- Row 34: Klein2003b, 6m ----
$ws.Range("A34").Value = "Klein2003b"

$ws.Range("B33").Value = 98.7
$ws.Range("C33").Value = -6.7
$ws.Range("D33").Value = 33.9
$ws.Range("F33").Value = 129.5
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 30
$ws.Range("I33").Value = 19.2
$ws.Range("L33").Value = 33
$ws.Range("M33").Value = "3m"
$ws.Range("P33").Value = "Cannot find actual macro intake, individuals instructed to eat less than 20 g CHO/day to start"
$ws.Range("Q33").Value = "x"
$ws.Range("S33").Value = $nejmLink

$ws.Range("B34").Value = 98.7
$ws.Range("C34").Value = -6.9
$ws.Range("D34").Value = 33.9
$ws.Range("F34").Value = 129.5
$ws.Range("G34").Value = 3.5
$ws.Range("H34").Value = 30
$ws.Range("I34").Value = 9.3000000000000007
$ws.Range("L34").Value = 33
$ws.Range("M34").Value = "6m"

# Reflect the final selection left after entering the data.
$null = $ws.Range("S39").Select()
